$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 9478
$ws.Range("C2").Value = 9447
$ws.Range("D2").Value = 8374
$ws.Range("E2").Value = 0.8864189689848629
$ws.Range("F2").Value = 0.8835197299008229
$ws.Range("G2").Value = 0.09681583147536164
$ws.Range("H2").Value = 0.0855386972752351
$ws.Range("I2").Value = 41174860.32985197
$ws.Range("J2").Value = 14390491.49810799
$ws.Range("L2").Value = 14390491.49810799
$ws.Range("M2").Value = 55565351.82795996
$ws.Range("N2").Value = 799515284.2372
$ws.Range("O2").Value = 781815477.2332001
$ws.Range("P2").Value = 0.01799901988345056
$ws.Range("Q2").Value = 0.01840650628845966
$ws.Range("B3").Value = 9666
$ws.Range("C3").Value = 9645
$ws.Range("D3").Value = 8563
$ws.Range("E3").Value = 0.887817522032141
$ws.Range("F3").Value = 0.8858886819780675
$ws.Range("G3").Value = 0.09539772302800158
$ws.Range("H3").Value = 0.08451176311698505
$ws.Range("I3").Value = 43125250.93561375
$ws.Range("J3").Value = 15095706.59566786
$ws.Range("L3").Value = 15095706.59566786
$ws.Range("M3").Value = 58220957.53128161
$ws.Range("N3").Value = 837887628.1232281
$ws.Range("O3").Value = 820407452.099158
$ws.Range("P3").Value = 0.01801638559752995
$ws.Range("Q3").Value = 0.01840025533293586
$ws.Range("B4").Value = 9858
$ws.Range("C4").Value = 9842
$ws.Range("D4").Value = 8755
$ws.Range("E4").Value = 0.889554968502337
$ws.Range("F4").Value = 0.8881111787380808
$ws.Range("G4").Value = 0.09419949378696643
$ws.Range("H4").Value = 0.08365962346367327
$ws.Range("I4").Value = 45180622.83051215
$ws.Range("J4").Value = 15803405.31809689
$ws.Range("L4").Value = 15803405.31809689
$ws.Range("M4").Value = 60984028.14860904
$ws.Range("N4").Value = 875120136.5024129
$ws.Range("O4").Value = 857671188.4964591
$ws.Range("P4").Value = 0.01805855523020903
$ws.Range("Q4").Value = 0.01842594869696049
$ws.Range("B5").Value = 10054
$ws.Range("C5").Value = 10024
$ws.Range("D5").Value = 8885
$ws.Range("E5").Value = 0.8863727055067837
$ws.Range("F5").Value = 0.8837278695046747
$ws.Range("G5").Value = 0.09334574158803877
$ws.Range("H5").Value = 0.08249223334093142
$ws.Range("I5").Value = 47130445.73264639
$ws.Range("J5").Value = 16445176.013788
$ws.Range("L5").Value = 16445176.013788
$ws.Range("M5").Value = 63575621.7464344
$ws.Range("N5").Value = 913371970.5379409
$ws.Range("O5").Value = 895885865.0749676
$ws.Range("P5").Value = 0.01800490549770476
$ws.Range("Q5").Value = 0.01835632936614294
$ws.Range("B6").Value = 10254
$ws.Range("C6").Value = 10235
$ws.Range("D6").Value = 9103
$ws.Range("E6").Value = 0.8893991206643869
$ws.Range("F6").Value = 0.8877511215135557
$ws.Range("G6").Value = 0.09209148262974684
$ws.Range("H6").Value = 0.08175431698640388
$ws.Range("I6").Value = 49481204.14574344
$ws.Range("J6").Value = 17249441.50190688
$ws.Range("L6").Value = 17249441.50190688
$ws.Range("M6").Value = 66730645.64765032
$ws.Range("N6").Value = 954978001.9095395
$ws.Range("O6").Value = 937386176.0363579
$ws.Range("P6").Value = 0.01806265847738432
$ws.Range("Q6").Value = 0.01840163845262194
